# Generate Report for Archive
#
# The localization status of e2e\bf804513-be9e-4e4a-bb07-552230a1ed0b.md
# moved from "Ready for handoff" to "In Translation". Update the Status
# column for that row on every sheet that tracks it:
#   - Overview sheet: zh-cn (E6) and de-de (F6) status columns
#   - zh-cn sheet: Status column (C6)
#   - de-de sheet: Status column (C6)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E6").Value = "In Translation"
$overview.Range("F6").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C6").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C6").Value = "In Translation"
